# Update "想去人数" (want-to-go count) figures in the "展览" and "全部类型"
# sheets to reflect the latest scrape (output generated at 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value  = 8104
$ws1.Range("F5").Value  = 104
$ws1.Range("F11").Value = 172
$ws1.Range("F13").Value = 463
$ws1.Range("F15").Value = 84
$ws1.Range("F17").Value = 6051
$ws1.Range("F19").Value = 290
$ws1.Range("F20").Value = 2151
$ws1.Range("F21").Value = 72
$ws1.Range("F22").Value = 120
$ws1.Range("F24").Value = 435

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value  = 8104
$ws4.Range("F5").Value  = 104
$ws4.Range("F12").Value = 172
$ws4.Range("F14").Value = 463
$ws4.Range("F16").Value = 84
$ws4.Range("F19").Value = 6051
$ws4.Range("F22").Value = 290
$ws4.Range("F23").Value = 2151
$ws4.Range("F24").Value = 72
$ws4.Range("F25").Value = 120
$ws4.Range("F27").Value = 435
